# Applies the text corrections described in the commit:
# "From 1.2.4 to 1.2.5 change and minor updates"
#
# The workbook stores its text as shared strings, so the same literal
# value can show up in several cells. Update every occurrence.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Version bump: "1.0" -> "1.2.5"
$ws.Range("D2").Value = "1.2.5"

# 2) Precondition text: fix accent + add trailing period
#    "O usuario ... sistema" -> "O usuário ... sistema."
$precondCells = @("B8", "B19", "B27", "B38", "B50")
foreach ($addr in $precondCells) {
    $ws.Range($addr).Value = "O usuário devidamente autenticado e na tela inicial do sistema."
}

# 3) "histório" -> "histórico" (typo fix)
$historicoCells = @("B12", "B31", "B42", "B54")
foreach ($addr in $historicoCells) {
    $ws.Range($addr).Value = "Chefe Verifica o histórico da tramitação da prestação de contas e clica para analisar a prestação de contas."
}

# 4) Add trailing period to "SYSTEM Exibe a tela para prestação de contas"
$telaPrestacaoCells = @("D12", "D31", "D42", "D54")
foreach ($addr in $telaPrestacaoCells) {
    $ws.Range($addr).Value = "SYSTEM Exibe a tela para prestação de contas."
}

# 5) Add trailing period to "SYSTEM Exibe a tela de Detalhar Diárias"
$ws.Range("D33").Value = "SYSTEM Exibe a tela de Detalhar Diárias."

# 6) Add semicolon before "Exibe mensagem de erro" in MSG203 sentence
$ws.Range("D56").Value = "SYSTEM Identifica que campos obrigatórios do parecer/análise não foram devidamente preenchidos; Exibe mensagem de erro (MSG203 - Campos obrigatórios) para o usuário."
